# Update cryptos list (Price / Volume(1h) columns) to the latest scraped
# values, matching the GitHub Actions refresh commit.
# For cells whose new text looks like a plain number (e.g. "0.990"),
# the Price column is stored as text, so we force the cell's number
# format to Text ("@") before assigning the value to stop Excel from
# auto-converting it to a numeric value, then restore the default
# "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.785.72'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '1.623.19'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.78%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.990'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.44'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.96%  '
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0909'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.853.75'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('D13').Value = '1.631.59'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.568'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.04%  '
$ws.Range('D16').Value = '29.781.93'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '9.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +18.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.39'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').Value = '  +2.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.68'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.56%  '
$ws.Range('E27').Value = '  +2.12%  '
$ws.Range('E28').Value = '  +3.58%  '
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('E32').Value = '  +3.14%  '
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('D34').Value = '1.424.46'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('E35').Value = '  +7.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0170'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.559'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0504'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.829'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '53.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +19.27%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.990'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('D49').Value = '1.763.31'
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '88.19'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('E51').Value = '  +3.05%  '
